$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H133").Value = 108999
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 108999
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 108999
$ws.Range("N133").Value = -119119

$ws.Range("H138").Value = 2554.7212
$ws.Range("I138").Value = 2024.5
$ws.Range("J138").Value = 2612.5637
$ws.Range("K138").Value = 6073.5
$ws.Range("L138").Value = 7837.6911
$ws.Range("M138").Value = -933.5
$ws.Range("N138").Value = -18117.6911

$ws.Range("H141").Value = 2976.9375
$ws.Range("I141").Value = 2835.0667
$ws.Range("J141").Value = 5105
$ws.Range("K141").Value = 8505.2001
$ws.Range("L141").Value = 15315
$ws.Range("M141").Value = -3325.2001
$ws.Range("N141").Value = -25675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2829.889
$ws.Range("I32").Value = 2825.0857
$ws.Range("J32").Value = 2998
$ws.Range("K32").Value = 2825.0857
$ws.Range("L32").Value = 2998
$ws.Range("M32").Value = -2538.0857
$ws.Range("N32").Value = -3572

$ws.Range("H110").Value = 4909.625
$ws.Range("I110").Value = 4212.6665
$ws.Range("J110").Value = 7000.5
$ws.Range("K110").Value = 4212.6665
$ws.Range("L110").Value = 7000.5
$ws.Range("M110").Value = -2167.6665
$ws.Range("N110").Value = -11090.5

$ws.Range("H122").Value = 3197.611
$ws.Range("I122").Value = 3039.8572
$ws.Range("J122").Value = 3749.75
$ws.Range("K122").Value = 9119.571599999999
$ws.Range("L122").Value = 11249.25
$ws.Range("M122").Value = -6669.571599999999
$ws.Range("N122").Value = -16149.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 90914560
$ws.Range("I94").Value = 133340664
$ws.Range("J94").Value = 1477.7142
$ws.Range("K94").Value = 133340664
$ws.Range("L94").Value = 1477.7142
$ws.Range("M94").Value = -133340213
$ws.Range("N94").Value = -2379.7142

$ws.Range("H99").Value = 58384.39
$ws.Range("I99").Value = 65088.75
$ws.Range("J99").Value = 4749.5
$ws.Range("K99").Value = 65088.75
$ws.Range("L99").Value = 4749.5
$ws.Range("M99").Value = -63590.75
$ws.Range("N99").Value = -7745.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 40000000
$ws.Range("I6").Value = 40000000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 40000000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -39999887
$ws.Range("N6").ClearContents()

$ws.Range("H31").Value = 4129.921
$ws.Range("I31").Value = 2881.5386
$ws.Range("J31").Value = 6834.75
$ws.Range("K31").Value = 2881.5386
$ws.Range("L31").Value = 6834.75
$ws.Range("M31").Value = -2586.5386
$ws.Range("N31").Value = -7424.75

$ws.Range("H34").Value = 4129.921
$ws.Range("I34").Value = 2881.5386
$ws.Range("J34").Value = 6834.75
$ws.Range("K34").Value = 2881.5386
$ws.Range("L34").Value = 6834.75
$ws.Range("M34").Value = -2679.5386
$ws.Range("N34").Value = -7238.75

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H74").Value = 88000
$ws.Range("I74").Value = 88000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 88000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -87126
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 88000
$ws.Range("I77").Value = 88000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 264000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -259632
$ws.Range("N77").ClearContents()

$ws.Range("H132").Value = 15157151
$ws.Range("I132").Value = 4979.4
$ws.Range("J132").Value = 27783962
$ws.Range("K132").Value = 14938.2
$ws.Range("L132").Value = 83351886
$ws.Range("M132").Value = -12408.2
$ws.Range("N132").Value = -83356946

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 165.48276
$ws.Range("I2").Value = 111.933334
$ws.Range("J2").Value = 222.85715
$ws.Range("K2").Value = 671.600004
$ws.Range("L2").Value = 1337.1429
$ws.Range("M2").Value = -558.600004
$ws.Range("N2").Value = -1563.1429

$ws.Range("H18").Value = 1152.4
$ws.Range("I18").Value = 1152.4
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 3457.2
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -3288.2
$ws.Range("N18").ClearContents()

$ws.Range("H60").Value = 2502323.5
$ws.Range("I60").Value = 6667133.5
$ws.Range("J60").Value = 3437.6
$ws.Range("K60").Value = 20001400.5
$ws.Range("L60").Value = 10312.8
$ws.Range("M60").Value = -20001149.5
$ws.Range("N60").Value = -10814.8

$ws.Range("H92").Value = 408.66666
$ws.Range("I92").Value = 262.5
$ws.Range("J92").Value = 437.9
$ws.Range("K92").Value = 787.5
$ws.Range("L92").Value = 1313.7
$ws.Range("M92").Value = 460.5
$ws.Range("N92").Value = -3809.7

$ws.Range("H113").Value = 1889.6666
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1889.6666
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5668.9998
$ws.Range("N113").Value = -10008.9998

$ws.Range("H122").Value = 1849.3846
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1849.3846
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 16644.4614
$ws.Range("N122").Value = -21544.4614

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H129").Value = 186460
$ws.Range("I129").Value = 1630
$ws.Range("J129").Value = 278875
$ws.Range("K129").Value = 4890
$ws.Range("L129").Value = 836625
$ws.Range("M129").Value = 110
$ws.Range("N129").Value = -846625

$ws.Range("H134").Value = 3513.6316
$ws.Range("I134").Value = 1484.5
$ws.Range("J134").Value = 4989.364
$ws.Range("K134").Value = 4453.5
$ws.Range("L134").Value = 14968.092
$ws.Range("M134").Value = 616.5
$ws.Range("N134").Value = -25108.092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H97").Value = 1311.5
$ws.Range("I97").Value = 1417.8889
$ws.Range("J97").Value = 354
$ws.Range("K97").Value = 1417.8889
$ws.Range("L97").Value = 354
$ws.Range("M97").Value = -921.8888999999999
$ws.Range("N97").Value = -1346

$ws.Range("H102").Value = 7755.294
$ws.Range("I102").Value = 406
$ws.Range("J102").Value = 8735.200000000001
$ws.Range("K102").Value = 406
$ws.Range("L102").Value = 8735.200000000001
$ws.Range("M102").Value = 1216
$ws.Range("N102").Value = -11979.2

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H122").Value = 3208631.2
$ws.Range("I122").Value = 4810576
$ws.Range("J122").Value = 4741.375
$ws.Range("K122").Value = 14431728
$ws.Range("L122").Value = 14224.125
$ws.Range("M122").Value = -14429278
$ws.Range("N122").Value = -19124.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 36000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 36000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 36000
$ws.Range("N20").Value = -36452
$ws.Range("M20").ClearContents()

$ws.Range("H82").Value = 881.9048
$ws.Range("I82").Value = 845
$ws.Range("J82").Value = 955.7143
$ws.Range("K82").Value = 845
$ws.Range("L82").Value = 955.7143
$ws.Range("M82").Value = -484
$ws.Range("N82").Value = -1677.7143

$ws.Range("H85").Value = 881.9048
$ws.Range("I85").Value = 845
$ws.Range("J85").Value = 955.7143
$ws.Range("K85").Value = 845
$ws.Range("L85").Value = 955.7143
$ws.Range("M85").Value = 403
$ws.Range("N85").Value = -3451.7143

$ws.Range("H132").Value = 5233.64
$ws.Range("I132").Value = 5055.7646
$ws.Range("J132").Value = 5611.625
$ws.Range("K132").Value = 15167.2938
$ws.Range("L132").Value = 16834.875
$ws.Range("M132").Value = -12637.2938
$ws.Range("N132").Value = -21894.875
